$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cluster id / count pairs in rows 2-6 with the new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 69

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 68

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 65

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 63

$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 50

# Remove the now-unused rows 7-11 and shift remaining cells up
$ws.Range("A7:B11").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
